$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 39000
$ws.Range("J123").Value = 39000
$ws.Range("L123").Value = 39000
$ws.Range("N123").Value = -48800
$ws.Range("H132").Value = 2567.5
$ws.Range("I132").Value = 2567.5
$ws.Range("K132").Value = 7702.5
$ws.Range("M132").Value = -5172.5
$ws.Range("H138").Value = 2646.32
$ws.Range("I138").Value = 1184.5358
$ws.Range("J138").Value = 3214.7917
$ws.Range("K138").Value = 3553.6074
$ws.Range("L138").Value = 9644.375100000001
$ws.Range("M138").Value = 1586.3926
$ws.Range("N138").Value = -19924.3751
$ws.Range("H141").Value = 5240
$ws.Range("I141").Value = 5400
$ws.Range("K141").Value = 16200
$ws.Range("M141").Value = -11020

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1403.5
$ws.Range("I45").Value = 1400
$ws.Range("J45").Value = 1404.6666
$ws.Range("K45").Value = 1400
$ws.Range("L45").Value = 1404.6666
$ws.Range("M45").Value = -1023
$ws.Range("N45").Value = -2158.6666
$ws.Range("H74").Value = 1257.4166
$ws.Range("I74").Value = 1208.091
$ws.Range("J74").Value = 1800
$ws.Range("K74").Value = 1208.091
$ws.Range("L74").Value = 1800
$ws.Range("M74").Value = -334.0909999999999
$ws.Range("N74").Value = -3548
$ws.Range("H77").Value = 1257.4166
$ws.Range("I77").Value = 1208.091
$ws.Range("J77").Value = 1800
$ws.Range("K77").Value = 6040.455
$ws.Range("L77").Value = 9000
$ws.Range("M77").Value = -1672.455
$ws.Range("N77").Value = -17736
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H122").Value = 1535.0667
$ws.Range("I122").Value = 896
$ws.Range("J122").Value = 2094.25
$ws.Range("K122").Value = 2688
$ws.Range("L122").Value = 6282.75
$ws.Range("M122").Value = -238
$ws.Range("N122").Value = -11182.75
$ws.Range("H132").Value = 1810.2084
$ws.Range("I132").Value = 1140.6
$ws.Range("J132").Value = 2926.2222
$ws.Range("K132").Value = 3421.8
$ws.Range("L132").Value = 8778.6666
$ws.Range("M132").Value = -891.7999999999997
$ws.Range("N132").Value = -13838.6666

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 1820.6666
$ws.Range("I33").Value = 1820.6666
$ws.Range("K33").Value = 1820.6666
$ws.Range("M33").Value = -1484.6666
$ws.Range("H36").Value = 8575.632
$ws.Range("I36").Value = 2000
$ws.Range("J36").Value = 8940.944
$ws.Range("K36").Value = 2000
$ws.Range("L36").Value = 8940.944
$ws.Range("M36").Value = -1466
$ws.Range("N36").Value = -10008.944
$ws.Range("H94").Value = 1295.12
$ws.Range("I94").Value = 1344.7
$ws.Range("J94").Value = 1096.8
$ws.Range("K94").Value = 1344.7
$ws.Range("L94").Value = 1096.8
$ws.Range("M94").Value = -893.7
$ws.Range("N94").Value = -1998.8
$ws.Range("H134").Value = 2642.16
$ws.Range("I134").Value = 2224.1304
$ws.Range("J134").Value = 7449.5
$ws.Range("K134").Value = 6672.3912
$ws.Range("L134").Value = 22348.5
$ws.Range("M134").Value = -4137.3912
$ws.Range("N134").Value = -27418.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2632.3333
$ws.Range("I99").Value = 2516.4119
$ws.Range("J99").Value = 3125
$ws.Range("K99").Value = 2516.4119
$ws.Range("L99").Value = 3125
$ws.Range("M99").Value = -1018.4119
$ws.Range("N99").Value = -6121
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H126").Value = 2632.3333
$ws.Range("I126").Value = 2516.4119
$ws.Range("J126").Value = 3125
$ws.Range("K126").Value = 7549.2357
$ws.Range("L126").Value = 9375
$ws.Range("M126").Value = -5079.2357
$ws.Range("N126").Value = -14315
$ws.Range("H134").Value = 35715820
$ws.Range("I134").Value = 1562.6
$ws.Range("J134").Value = 125001464
$ws.Range("K134").Value = 4687.799999999999
$ws.Range("L134").Value = 375004392
$ws.Range("M134").Value = -2152.799999999999
$ws.Range("N134").Value = -375009462

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2959.0527
$ws.Range("I126").Value = 2697.2144
$ws.Range("J126").Value = 3692.2
$ws.Range("K126").Value = 8091.6432
$ws.Range("L126").Value = 11076.6
$ws.Range("M126").Value = -5621.6432
$ws.Range("N126").Value = -16016.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2154.889
$ws.Range("I7").Value = 2291.5
$ws.Range("J7").Value = 1881.6666
$ws.Range("K7").Value = 2291.5
$ws.Range("L7").Value = 1881.6666
$ws.Range("M7").Value = -2179.5
$ws.Range("N7").Value = -2105.6666
$ws.Range("H40").Value = 10101010
$ws.Range("I40").Value = 10101010
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 10101010
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -10100874
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 2057.4285
$ws.Range("I46").Value = 2233.3333
$ws.Range("K46").Value = 2233.3333
$ws.Range("M46").Value = -2045.3333
$ws.Range("H55").Value = 298.0909
$ws.Range("I55").Value = 248.14285
$ws.Range("J55").Value = 385.5
$ws.Range("K55").Value = 248.14285
$ws.Range("L55").Value = 385.5
$ws.Range("M55").Value = -75.14285000000001
$ws.Range("N55").Value = -731.5
$ws.Range("H93").Value = 2928.7144
$ws.Range("I93").Value = 2140.6
$ws.Range("J93").Value = 3366.5557
$ws.Range("K93").Value = 2140.6
$ws.Range("L93").Value = 3366.5557
$ws.Range("M93").Value = -892.5999999999999
$ws.Range("N93").Value = -5862.5557
$ws.Range("H100").Value = 16033631
$ws.Range("I100").Value = 18705436
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 18705436
$ws.Range("L100").Value = 2800
$ws.Range("M100").Value = -18704895
$ws.Range("N100").Value = -3882
$ws.Range("H122").Value = 2726.25
$ws.Range("I122").Value = 1950
$ws.Range("J122").Value = 3502.5
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 10507.5
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -15407.5
$ws.Range("H126").Value = 2154.889
$ws.Range("I126").Value = 2291.5
$ws.Range("J126").Value = 1881.6666
$ws.Range("K126").Value = 6874.5
$ws.Range("L126").Value = 5644.9998
$ws.Range("M126").Value = -4404.5
$ws.Range("N126").Value = -10584.9998
$ws.Range("H132").Value = 4118.4546
$ws.Range("I132").Value = 2500
$ws.Range("J132").Value = 4280.3
$ws.Range("K132").Value = 7500
$ws.Range("L132").Value = 12840.9
$ws.Range("M132").Value = -4970
$ws.Range("N132").Value = -17900.9

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 624.2143
$ws.Range("I126").Value = 643
$ws.Range("J126").Value = 380
$ws.Range("K126").Value = 1929
$ws.Range("L126").Value = 1140
$ws.Range("M126").Value = 541
$ws.Range("N126").Value = -6080
